$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.00391733333333
$ws.Range("H2").Value = 105.011752
$ws.Range("I2").Value = 0.9591895364534718
$ws.Range("J2").Value = 0.9591895364534718
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 70.95697755313243
$ws.Range("R2").Value = 638.612797978192
$ws.Range("S2").Value = 0.006327087134619072
$ws.Range("T2").Value = 0.006327087134619072
$ws.Range("G3").Value = 35.00391733333333
$ws.Range("H3").Value = 105.011752
$ws.Range("I3").Value = 0.9591895364534718
$ws.Range("J3").Value = 0.9591895364534718
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 8976.555264492092
$ws.Range("R3").Value = 80788.99738042883
$ws.Range("S3").Value = 0.8004208928521047
$ws.Range("T3").Value = 0.8004208928521046
$ws.Range("G4").Value = 35.00391733333333
$ws.Range("H4").Value = 105.011752
$ws.Range("I4").Value = 0.9591895364534718
$ws.Range("J4").Value = 0.9591895364534718
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 1709.600621934039
$ws.Range("R4").Value = 15386.40559740635
$ws.Range("S4").Value = 0.1524415564667482
$ws.Range("T4").Value = 0.1524415564667481
$ws.Range("I5").Value = 0.0008369499257158872
$ws.Range("J5").Value = 0.0008369499257158872
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 0.061914183626
$ws.Range("R5").Value = 0.557227652634
$ws.Range("S5").Value = 0.000005520759877027964
$ws.Range("T5").Value = 0.000005520759877027963
$ws.Range("I6").Value = 0.0008369499257158872
$ws.Range("J6").Value = 0.0008369499257158872
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.000698414840189939
$ws.Range("T6").Value = 0.0006984148401899388
$ws.Range("I7").Value = 0.0008369499257158872
$ws.Range("J7").Value = 0.0008369499257158872
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 1.491728234257
$ws.Range("R7").Value = 13.425554108313
$ws.Range("S7").Value = 0.0001330143256489204
$ws.Range("T7").Value = 0.0001330143256489204
$ws.Range("G8").Value = 1.458762333333333
$ws.Range("H8").Value = 4.376287
$ws.Range("I8").Value = 0.03997351362081222
$ws.Range("J8").Value = 0.03997351362081222
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 2.957079493589111
$ws.Range("R8").Value = 26.613715442302
$ws.Range("S8").Value = 0.0002636766709225144
$ws.Range("T8").Value = 0.0002636766709225144
$ws.Range("G9").Value = 1.458762333333333
$ws.Range("H9").Value = 4.376287
$ws.Range("I9").Value = 0.03997351362081222
$ws.Range("J9").Value = 0.03997351362081222
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 374.0912932180991
$ws.Range("R9").Value = 3366.821638962891
$ws.Range("S9").Value = 0.03335694797204277
$ws.Range("T9").Value = 0.03335694797204276
$ws.Range("G10").Value = 1.458762333333333
$ws.Range("H10").Value = 4.376287
$ws.Range("I10").Value = 0.03997351362081222
$ws.Range("J10").Value = 0.03997351362081222
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 71.24633990452656
$ws.Range("R10").Value = 641.217059140739
$ws.Range("S10").Value = 0.006352888977846935
$ws.Range("T10").Value = 0.006352888977846933
